# Corrections to two PowerPoint slides (slide 13 and slide 14):
#   - Slide 13: "  BLE L2" -> "  BG L2"
#   - Slide 14: two runs ("  " + "BG L2") -> single run "  BLE L2"

$p = $ppt.ActivePresentation

# --- Slide 13: fix "  BLE L2" run to "  BG L2" -----------------------------
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)          # "Content Placeholder 9"
$tr13 = $sh13.TextFrame.TextRange

$len13 = $tr13.Length
$needle13 = "  BLE L2"
$nlen13 = $needle13.Length
$pos13 = -1
for ($i = 1; $i -le ($len13 - $nlen13 + 1); $i++) {
    $probe = $tr13.Characters($i, $nlen13)
    if ($probe.Text -eq $needle13) {
        $pos13 = $i
        break
    }
}

if ($pos13 -gt 0) {
    $run13 = $tr13.Characters($pos13, $nlen13)
    $run13.Text = "  BG L2"
}

# --- Slide 14: merge "  " + "BG L2" runs into single "  BLE L2" run -------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(2)          # "Content Placeholder 9"
$tr14 = $sh14.TextFrame.TextRange

$len14 = $tr14.Length
$needle14 = "  BG L2"
$nlen14 = $needle14.Length
$pos14 = -1
for ($i = 1; $i -le ($len14 - $nlen14 + 1); $i++) {
    $probe = $tr14.Characters($i, $nlen14)
    if ($probe.Text -eq $needle14) {
        $pos14 = $i
        break
    }
}

if ($pos14 -gt 0) {
    # Clear the leading "  " run (2 spaces) entirely so it disappears,
    # leaving the "BG L2" run (and its original, dirty-free formatting)
    # as the sole remaining run.
    $spaces = $tr14.Characters($pos14, 2)
    $spaces.Text = ""

    # "BG L2" now starts where "  " used to be (position $pos14).
    $word = $tr14.Characters($pos14, 5)
    $word.Text = "  BLE L2"
}
